$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 49843.07
$ws.Range("J17").Value = 49843.07
$ws.Range("L17").Value = 149529.21
$ws.Range("N17").Value = -149865.21

$ws.Range("H98").Value = 1671.6897
$ws.Range("I98").Value = 1116.1904
$ws.Range("J98").Value = 3129.875
$ws.Range("K98").Value = 1116.1904
$ws.Range("L98").Value = 3129.875
$ws.Range("M98").Value = 381.8096
$ws.Range("N98").Value = -6125.875

$ws.Range("H122").Value = 1671.6897
$ws.Range("I122").Value = 1116.1904
$ws.Range("J122").Value = 3129.875
$ws.Range("K122").Value = 3348.5712
$ws.Range("L122").Value = 9389.625
$ws.Range("M122").Value = -898.5711999999999
$ws.Range("N122").Value = -14289.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1662.3784
$ws.Range("I45").Value = 1060.7142
$ws.Range("J45").Value = 3534.2222
$ws.Range("K45").Value = 1060.7142
$ws.Range("L45").Value = 3534.2222
$ws.Range("M45").Value = -683.7141999999999
$ws.Range("N45").Value = -4288.2222

$ws.Range("H122").Value = 3302.2173
$ws.Range("I122").Value = 2336.2222
$ws.Range("J122").Value = 6779.8
$ws.Range("K122").Value = 7008.6666
$ws.Range("L122").Value = 20339.4
$ws.Range("M122").Value = -4558.6666
$ws.Range("N122").Value = -25239.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1637.3043
$ws.Range("I105").Value = 1407.0588
$ws.Range("J105").Value = 2289.6667
$ws.Range("K105").Value = 1407.0588
$ws.Range("L105").Value = 2289.6667
$ws.Range("M105").Value = 339.9412
$ws.Range("N105").Value = -5783.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1351.2858
$ws.Range("I16").Value = 1086.1818
$ws.Range("J16").Value = 1522.8235
$ws.Range("K16").Value = 1086.1818
$ws.Range("L16").Value = 1522.8235
$ws.Range("M16").Value = -799.1818000000001
$ws.Range("N16").Value = -2096.8235

$ws.Range("H31").Value = 1786.5471
$ws.Range("I31").Value = 1159.5714
$ws.Range("J31").Value = 3005.6667
$ws.Range("K31").Value = 1159.5714
$ws.Range("L31").Value = 3005.6667
$ws.Range("M31").Value = -864.5714
$ws.Range("N31").Value = -3595.6667

$ws.Range("H34").Value = 1786.5471
$ws.Range("I34").Value = 1159.5714
$ws.Range("J34").Value = 3005.6667
$ws.Range("K34").Value = 1159.5714
$ws.Range("L34").Value = 3005.6667
$ws.Range("M34").Value = -957.5714
$ws.Range("N34").Value = -3409.6667

$ws.Range("H99").Value = 1861.9615
$ws.Range("I99").Value = 1123.4445
$ws.Range("J99").Value = 2252.9412
$ws.Range("K99").Value = 1123.4445
$ws.Range("L99").Value = 2252.9412
$ws.Range("M99").Value = 374.5554999999999
$ws.Range("N99").Value = -5248.9412

$ws.Range("H113").Value = 1351.2858
$ws.Range("I113").Value = 1086.1818
$ws.Range("J113").Value = 1522.8235
$ws.Range("K113").Value = 1086.1818
$ws.Range("L113").Value = 1522.8235
$ws.Range("M113").Value = 1083.8182
$ws.Range("N113").Value = -5862.8235

$ws.Range("H122").Value = 2340.8572
$ws.Range("I122").Value = 1264.6923
$ws.Range("J122").Value = 4089.625
$ws.Range("K122").Value = 3794.0769
$ws.Range("L122").Value = 12268.875
$ws.Range("M122").Value = -1344.0769
$ws.Range("N122").Value = -17168.875

$ws.Range("H126").Value = 1861.9615
$ws.Range("I126").Value = 1123.4445
$ws.Range("J126").Value = 2252.9412
$ws.Range("K126").Value = 3370.3335
$ws.Range("L126").Value = 6758.823600000001
$ws.Range("M126").Value = -900.3335000000002
$ws.Range("N126").Value = -11698.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2307.45
$ws.Range("I131").Value = 2720.75
$ws.Range("J131").Value = 1687.5
$ws.Range("K131").Value = 8162.25
$ws.Range("L131").Value = 5062.5
$ws.Range("M131").Value = -3122.25
$ws.Range("N131").Value = -15142.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 892.5333000000001
$ws.Range("I107").Value = 201.25
$ws.Range("K107").Value = 201.25
$ws.Range("M107").Value = 1718.75

$ws.Range("H113").Value = 3188.875
$ws.Range("J113").Value = 4325
$ws.Range("L113").Value = 4325
$ws.Range("N113").Value = -8665

$ws.Range("H122").Value = 4904.0713
$ws.Range("I122").Value = 3107
$ws.Range("J122").Value = 6701.143
$ws.Range("K122").Value = 9321
$ws.Range("L122").Value = 20103.429
$ws.Range("M122").Value = -6871
$ws.Range("N122").Value = -25003.429

$ws.Range("H126").Value = 4910.5454
$ws.Range("I126").Value = 4660.4443
$ws.Range("J126").Value = 5210.6665
$ws.Range("K126").Value = 13981.3329
$ws.Range("L126").Value = 15631.9995
$ws.Range("M126").Value = -11511.3329
$ws.Range("N126").Value = -20571.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2200.7646
$ws.Range("I7").Value = 1678.8
$ws.Range("J7").Value = 2418.25
$ws.Range("K7").Value = 1678.8
$ws.Range("L7").Value = 2418.25
$ws.Range("M7").Value = -1566.8
$ws.Range("N7").Value = -2642.25

$ws.Range("H40").Value = 3283.3333
$ws.Range("I40").Value = 2033.3334
$ws.Range("J40").Value = 4533.3335
$ws.Range("K40").Value = 2033.3334
$ws.Range("L40").Value = 4533.3335
$ws.Range("M40").Value = -1897.3334
$ws.Range("N40").Value = -4805.3335

$ws.Range("H61").Value = 4689.4443
$ws.Range("I61").Value = 4017.3333
$ws.Range("K61").Value = 4017.3333
$ws.Range("M61").Value = -3815.3333

$ws.Range("H113").Value = 4689.4443
$ws.Range("I113").Value = 4017.3333
$ws.Range("K113").Value = 4017.3333
$ws.Range("M113").Value = -1847.3333

$ws.Range("H122").Value = 2718.0908
$ws.Range("I122").Value = 2189.9
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 6569.700000000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -4119.700000000001
$ws.Range("N122").Value = -28900

$ws.Range("H126").Value = 2200.7646
$ws.Range("I126").Value = 1678.8
$ws.Range("J126").Value = 2418.25
$ws.Range("K126").Value = 5036.4
$ws.Range("L126").Value = 7254.75
$ws.Range("M126").Value = -2566.4
$ws.Range("N126").Value = -12194.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2075.1428
$ws.Range("I122").Value = 1585.875
$ws.Range("J122").Value = 2727.5
$ws.Range("K122").Value = 4757.625
$ws.Range("L122").Value = 8182.5
$ws.Range("M122").Value = -2307.625
$ws.Range("N122").Value = -13082.5

$ws.Range("H126").Value = 2342.5356
$ws.Range("I126").Value = 2046.6923
$ws.Range("J126").Value = 2598.9333
$ws.Range("K126").Value = 6140.0769
$ws.Range("L126").Value = 7796.7999
$ws.Range("M126").Value = -3670.0769
$ws.Range("N126").Value = -12736.7999
